$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.009621207843739724, 0.2335480718435533, 0.1699024050833353),
    @(0.01924044167894088, 0.2316540184638475, 0.1679318969715666),
    @(0.01924044167894088, 0.2316540184638475, 0.1679318969715666)
)

$startRow = 80
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
